$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-10, columns E (5) through T (20).
# Column layout:
# A=1 Sending cluster, B=2 Ligand symbol, C=3 Receptor symbol, D=4 Target cluster,
# E=5 Ligand-expressing cells, F=6 Ligand detection rate,
# G=7 Ligand average expression value, H=8 Ligand total expression value,
# I=9 Ligand derived specificity avg, J=10 Ligand derived specificity total,
# K=11 Receptor-expressing cells, L=12 Receptor detection rate,
# M=13 Receptor average expression value, N=14 Receptor total expression value,
# O=15 Receptor derived specificity avg, P=16 Receptor derived specificity total,
# Q=17 Edge average expression weight, R=18 Edge total expression weight,
# S=19 Edge average expression derived specificity, T=20 Edge total expression derived specificity

$data = @{
    2  = @(3, 1, 0.8164013333333333, 2.449204, 0.05618115571687973, 0.05618115571687973, 2, 0.6666666666666666, 0.05968133333333333, 0.179044, 0.02602747651633847, 0.02602747651633848, 0.04872392010844444, 0.438515280976, 0.001462253711081842, 0.001462253711081842)
    3  = @(3, 1, 0.8164013333333333, 2.449204, 0.05618115571687973, 0.05618115571687973, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.144012433133819, 0.144012433133819, 0.2695939532297777, 2.426345579068, 0.008090784931057816, 0.008090784931057816)
    4  = @(3, 1, 0.8164013333333333, 2.449204, 0.05618115571687973, 0.05618115571687973, 3, 1, 1.903109, 5.709327, 0.8299600903498424, 0.8299600903498425, 1.553700725078667, 13.983306525708, 0.04662811707474007, 0.04662811707474007)
    5  = @(3, 1, 12.878362, 38.635086, 0.8862323361798529, 0.8862323361798529, 2, 0.6666666666666666, 0.05968133333333333, 0.179044, 0.02602747651633847, 0.02602747651633848, 0.7685978153093334, 6.917380337784, 0.0230663913179409, 0.0230663913179409)
    6  = @(3, 1, 12.878362, 38.635086, 0.8862323361798529, 0.8862323361798529, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.144012433133819, 0.144012433133819, 4.252722749151333, 38.274504742362, 0.1276284750551293, 0.1276284750551293)
    7  = @(3, 1, 12.878362, 38.635086, 0.8862323361798529, 0.8862323361798529, 3, 1, 1.903109, 5.709327, 0.8299600903498424, 0.8299600903498425, 24.508926627458, 220.580339647122, 0.7355374698067826, 0.7355374698067827)
    8  = @(3, 1, 0.8368233333333334, 2.51047, 0.05758650810326746, 0.05758650810326746, 2, 0.6666666666666666, 0.05968133333333333, 0.179044, 0.02602747651633847, 0.02602747651633848, 0.04994273229777778, 0.44948459068, 0.001498831487315729, 0.001498831487315729)
    9  = @(3, 1, 0.8368233333333334, 2.51047, 0.05758650810326746, 0.05758650810326746, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.144012433133819, 0.144012433133819, 0.2763377537211111, 2.48703978349, 0.008293173147631932, 0.008293173147631932)
    10 = @(3, 1, 0.8368233333333334, 2.51047, 0.05758650810326746, 0.05758650810326746, 3, 1, 1.903109, 5.709327, 0.8299600903498424, 0.8299600903498425, 1.592566017076667, 14.33309415369, 0.04779450346831979, 0.0477945034683198)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i  # starts at column E = 5
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}
